# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps for the
# f7f85e1d-... row (row 3) on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E3").Value = "2016-03-23 16:50:31"
$zhcn.Range("H3").Value = "2016-03-23 16:51:00"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E3").Value = "2016-03-23 16:50:36"
$dede.Range("H3").Value = "2016-03-23 16:51:08"
